$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (changed) date column C for rows 2-11 from 45202 to 45203
for ($row = 2; $row -le 11; $row++) {
    $ws.Cells.Item($row, 3).Value = 45203
}
